# Update the "Chart" (GSC export) data: drop the two oldest days
# (2025-09-28, 2025-09-29), shift every remaining day up, and append the
# newly observed day (2025-12-28) at the end of the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$lastOldRow = 92
$firstDataRow = 2
$lastNewRow = 91

# New date sequence for rows 2..91 (2025-09-30 .. 2025-12-28)
$newDates = @("2025-09-30","2025-10-01","2025-10-02","2025-10-03","2025-10-04","2025-10-05","2025-10-06","2025-10-07","2025-10-08","2025-10-09","2025-10-10","2025-10-11","2025-10-12","2025-10-13","2025-10-14","2025-10-15","2025-10-16","2025-10-17","2025-10-18","2025-10-19","2025-10-20","2025-10-21","2025-10-22","2025-10-23","2025-10-24","2025-10-25","2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18","2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27","2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06","2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11","2025-12-12","2025-12-13","2025-12-14","2025-12-15","2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23","2025-12-24","2025-12-25","2025-12-26","2025-12-27","2025-12-28")

# Capture the existing "HTTPS URLs" counts (column C) for rows 2..92 before
# overwriting anything, so the shift can be computed from the original data.
$oldCounts = @()
for ($r = $firstDataRow; $r -le $lastOldRow; $r++) {
    $oldCounts += $ws.Cells.Item($r, 3).Value()
}

# Row N (2..90) takes the count that used to live two rows further down;
# the brand-new final row (91, for 2025-12-28) has no data yet, so it is 0.
for ($r = $firstDataRow; $r -le $lastNewRow; $r++) {
    $dateText = $newDates[$r - $firstDataRow]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $dateText
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 2).Value = 0

    $offset = $r + 2 - $firstDataRow
    if ($offset -lt $oldCounts.Count) {
        $ws.Cells.Item($r, 3).Value = $oldCounts[$offset]
    } else {
        $ws.Cells.Item($r, 3).Value = 0
    }
}

# The series is now one row shorter; drop the old trailing row 92.
$ws.Range("A92:C92").ClearContents()
